# Update cryptos list prices / 1h-volume percentages (refresh snapshot).
# Price cells that look like a parseable decimal number (e.g. "141.34")
# are forced to text format first so they stay literal strings (matching
# the "957.30"-style / "1.00"-style display text) instead of being
# auto-coerced into real numbers by Excel's normal text-to-number input
# parsing. Price cells that aren't valid numbers anyway (multiple dots,
# the subscript PEPE price, etc.) don't need that treatment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.867.69"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.114.39"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.90"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.34"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.113.66"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.110"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").Value = "3.648.36"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.30"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "57.956.59"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "3.111.00"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.89"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.09"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "337.00"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.68"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "0.0₃0933"
$ws.Range("E28").Value = "  +4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.59"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.25"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.98"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.44"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "3.153.94"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.686"
$ws.Range("E42").Value = "  +5.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.92"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.96"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  +9.04%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "2.300.92"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +8.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.96"
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.02"
$ws.Range("E51").Value = "  +3.06%  "
